{"js": "// Attestation update \u2014 mirrors the author's diff:\n//  1. Rewrite the \"Je soussign\u00e9(e)...\" paragraph into the new\n//     \"Nous soussign\u00e9s, QuantFactory, ...\" wording (now also folding in the\n//     CDI / hiring-date sentence that used to live in the next paragraph).\n//  2. Rewrite the old \"Monsieur Med a \u00e9t\u00e9 embauch\u00e9(e)...\" paragraph into the\n//     shorter \"\u00c0 ce jour, Monsieur **Mouad Med** exerce toujours...\" line.\n//  3. Expand the closing date line into a 3-line sign-off block (date /\n//     \"Pour QuantFactory,\" / \"[Signature]\").\n//  4. Collapse the trailing \"Responsable des Ressources Humaines /\n//     QuantFactory / [Coordonn\u00e9es...]\" lines into a single\n//     \"**Service des Ressources Humaines**\" line.\n//\n// The whole letter lives in a single paragraph/run, with line breaks encoded\n// as <w:br/> between <w:t> runs of text; Word surfaces those breaks as the\n// vertical-tab character (\\u000b) in range.text / search() / insertText().\n\nconst body = context.document.body;\n\nasync function replaceText(oldText, newText, label) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find expected text for: \" + label);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- 1. Intro paragraph ---\nawait replaceText(\n  \"Je soussign\u00e9(e), [Nom du responsable RH], en ma qualit\u00e9 de repr\u00e9sentant(e) des Ressources Humaines chez QuantFactory, atteste par la pr\u00e9sente que Monsieur **Mouad Med** est employ\u00e9(e) au sein de notre entreprise en qualit\u00e9 de **Full-stack developer**, rattach\u00e9(e) au d\u00e9partement **IT**.  \",\n  \"Nous soussign\u00e9s, QuantFactory, attestons que Monsieur **Mouad Med** occupe au sein de notre entreprise le poste de **Full-stack developer** au sein du d\u00e9partement **IT**, et ce depuis son embauche en **Contrat \u00e0 Dur\u00e9e Ind\u00e9termin\u00e9e (CDI)** le **21 mai 2023**.  \",\n  \"intro paragraph\"\n);\n\n// --- 2. Tenure paragraph ---\nawait replaceText(\n  \"Monsieur Med a \u00e9t\u00e9 embauch\u00e9(e) en **Contrat \u00e0 Dur\u00e9e Ind\u00e9termin\u00e9e (CDI)** \u00e0 compter du **21 mai 2023** et occupe toujours ce poste \u00e0 ce jour.  \",\n  \"\u00c0 ce jour, Monsieur **Mouad Med** exerce toujours ses fonctions avec professionnalisme et d\u00e9vouement.  \",\n  \"tenure paragraph\"\n);\n\n// --- 3. Date line -> Date + \"Pour QuantFactory,\" + \"[Signature]\" ---\nawait replaceText(\n  \"Le [Date de r\u00e9daction compl\u00e8te, ex. 23 mai 2025]  \",\n  \"**Le [27/05/2025]**  \\u000b**Pour QuantFactory,**  \\u000b[Signature]  \",\n  \"date line\"\n);\n\n// --- 4. Collapse the last three signature lines into one ---\nawait replaceText(\n  \"Responsable des Ressources Humaines  \\u000bQuantFactory  \\u000b[Coordonn\u00e9es de l\u2019entreprise si n\u00e9cessaire]\",\n  \"**Service des Ressources Humaines**\",\n  \"footer lines\"\n);\n", "ps1": "# Attestation update \u2014 mirrors the author's diff:\n#  1. Rewrite the \"Je soussign\u00e9(e)...\" paragraph into the new\n#     \"Nous soussign\u00e9s, QuantFactory, ...\" wording (now also folding in the\n#     CDI / hiring-date sentence that used to live in the next paragraph).\n#  2. Rewrite the old \"Monsieur Med a \u00e9t\u00e9 embauch\u00e9(e)...\" paragraph into the\n#     shorter \"\u00c0 ce jour, Monsieur **Mouad Med** exerce toujours...\" line.\n#  3. Expand the closing date line into a 3-line sign-off block (date /\n#     \"Pour QuantFactory,\" / \"[Signature]\").\n#  4. Collapse the trailing \"Responsable des Ressources Humaines /\n#     QuantFactory / [Coordonn\u00e9es...]\" lines into a single\n#     \"**Service des Ressources Humaines**\" line.\n#\n# The whole letter lives in a single paragraph/run, with line breaks encoded\n# as <w:br/> between text runs; Word/COM surfaces those breaks as the\n# vertical-tab character (Chr 11) in Range.Text, and writing Chr 11 back into\n# Range.Text re-creates a <w:br/>.\n\n$d = $word.ActiveDocument\n$VT = [char]11\n\nfunction Replace-AttestationText($oldText, $newText, $label) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.Text = $oldText\n    $find.MatchWildcards = $false\n    $find.Forward = $true\n    $ok = $find.Execute()\n    if (-not $ok) {\n        throw \"Could not find expected text for: $label\"\n    }\n    $range.Text = $newText\n}\n\n# --- 1. Intro paragraph ---\nReplace-AttestationText \"Je soussign\u00e9(e), [Nom du responsable RH], en ma qualit\u00e9 de repr\u00e9sentant(e) des Ressources Humaines chez QuantFactory, atteste par la pr\u00e9sente que Monsieur **Mouad Med** est employ\u00e9(e) au sein de notre entreprise en qualit\u00e9 de **Full-stack developer**, rattach\u00e9(e) au d\u00e9partement **IT**.  \" \"Nous soussign\u00e9s, QuantFactory, attestons que Monsieur **Mouad Med** occupe au sein de notre entreprise le poste de **Full-stack developer** au sein du d\u00e9partement **IT**, et ce depuis son embauche en **Contrat \u00e0 Dur\u00e9e Ind\u00e9termin\u00e9e (CDI)** le **21 mai 2023**.  \" \"intro paragraph\"\n\n# --- 2. Tenure paragraph ---\nReplace-AttestationText \"Monsieur Med a \u00e9t\u00e9 embauch\u00e9(e) en **Contrat \u00e0 Dur\u00e9e Ind\u00e9termin\u00e9e (CDI)** \u00e0 compter du **21 mai 2023** et occupe toujours ce poste \u00e0 ce jour.  \" \"\u00c0 ce jour, Monsieur **Mouad Med** exerce toujours ses fonctions avec professionnalisme et d\u00e9vouement.  \" \"tenure paragraph\"\n\n# --- 3. Date line -> Date + \"Pour QuantFactory,\" + \"[Signature]\" ---\nReplace-AttestationText \"Le [Date de r\u00e9daction compl\u00e8te, ex. 23 mai 2025]  \" (\"**Le [27/05/2025]**  \" + $VT + \"**Pour QuantFactory,**  \" + $VT + \"[Signature]  \") \"date line\"\n\n# --- 4. Collapse the last three signature lines into one ---\nReplace-AttestationText (\"Responsable des Ressources Humaines  \" + $VT + \"QuantFactory  \" + $VT + \"[Coordonn\u00e9es de l\u2019entreprise si n\u00e9cessaire]\") \"**Service des Ressources Humaines**\" \"footer lines\"\n"}
